# Generate Report for Handback
# Refresh the handoff/handback timestamps for the file that was just
# processed (462cda04-9234-4fcc-8834-1276f7d60d8d, row 2 on every sheet),
# and roll the newest of those timestamps up into the Overview sheet's
# "Latest HO Xliff Generate Date" column.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-21 00:54:22"
$zhcn.Range("K2").Value = "2016-08-21 00:54:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-21 00:54:26"
$dede.Range("K2").Value = "2016-08-21 00:54:44"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-21 00:54:26"
